$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 171
$ws.Range("I2").Value = 434
$ws.Range("J2").Value = 1903
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 541
$ws.Range("M2").Value = 35
$ws.Range("N2").Value = 310
$ws.Range("P2").Value = 6
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 33
$ws.Range("S2").Value = 210
$ws.Range("T2").Value = 334
$ws.Range("U2").Value = 28
$ws.Range("V2").Value = 2957
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 3026
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 59
$ws.Range("AA2").Value = 21
